$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price/volume data from the Feb 16 2023 21:53:58 UTC GitHub Actions refresh
$updates = @{
    "D2" = "313.06"
    "E2" = "0.26%"
    "D3" = "48.58"
    "E3" = "8.99%"
    "E4" = "2.99%"
    "D5" = "0.07874"
    "E5" = "-1.72%"
    "D6" = "4.580"
    "E6" = "1.62%"
    "D7" = "1.321"
    "E7" = "21.71%"
    "D8" = "1.616"
    "E8" = "-1.87%"
    "D9" = "0.1244"
    "E9" = "-3.93%"
    "D10" = "0.1946"
    "E10" = "2.39%"
    "D11" = "0.09431"
    "E11" = "0.43%"
    "D12" = "0.04542"
    "E12" = "7.25%"
    "D13" = "0.1048"
    "E13" = "1.05%"
    "D14" = "0.001300"
    "E14" = "-0.59%"
    "D15" = "0.04197"
    "E15" = "-0.07%"
    "D16" = "0.005877"
    "E16" = "-0.86%"
    "E17" = "-1.00%"
    "E18" = "2.75%"
    "E19" = "2.25%"
    "D20" = "8.100"
    "E20" = "0.04%"
    "D21" = "0.1368"
    "E21" = "-0.17%"
    "E22" = "-2.22%"
    "D23" = "0.001291"
    "E23" = "1.31%"
    "D24" = "0.004156"
    "E24" = "-9.35%"
    "D25" = "0.0001353"
    "E25" = "0.82%"
    "D26" = "0.0003555"
    "D38" = "0.02623"
    "E38" = "-1.04%"
    "D39" = "0.05883"
    "E39" = "8.91%"
    "D40" = "0.01078"
    "E40" = "91.36%"
    "D41" = "0.008022"
    "E41" = "3.58%"
    "E42" = "2.00%"
    "D43" = "0.008244"
    "E43" = "12.31%"
    "D44" = "0.008555"
    "E44" = "8.62%"
    "D45" = "0.3137"
    "E45" = "0.72%"
    "D46" = "0.00006954"
    "E46" = "2.88%"
    "D47" = "0.00000000752"
    "E47" = "0.84%"
    "D48" = "0.05505"
    "E48" = "-3.27%"
    "D49" = "0.004017"
    "E49" = "1.08%"
    "D50" = "0.00002105"
    "E50" = "0.84%"
    "D51" = "0.0002004"
    "E51" = "0.84%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}

Write-Output ("Updated {0} cells" -f $updates.Count)
